# Add data for 2021-10-23
# Updates the "through" date from 10-14 to 10-15, refreshes the October
# row with new counts, and rolls the new deltas into the Total row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename sheet (tab name + workbook.xml sheet name)
$ws.Name = "Through 2021-10-15"

# 2. Update the row label for October
$ws.Range("A11").Value = "October (through 10-15)"

# 3. Update October row (row 11) values for 2016-2021 (columns C-H)
$ws.Range("C11").Value = 26
$ws.Range("D11").Value = 29
$ws.Range("E11").Value = 38
$ws.Range("F11").Value = 19
$ws.Range("G11").Value = 75
$ws.Range("H11").Value = 95

# 4. Update Total row (row 12) values for 2016-2021 (columns C-H)
$ws.Range("C12").Value = 455
$ws.Range("D12").Value = 656
$ws.Range("E12").Value = 586
$ws.Range("F12").Value = 441
$ws.Range("G12").Value = 976
$ws.Range("H12").Value = 1345
